$wb = $excel.ActiveWorkbook

# --- subscribers sheet: add the missing Paramount Plus 2020 figure ---
$subs = $wb.Worksheets.Item("subscribers")
$subs.Range("E9").Value = 8

# --- notes sheet: populate with headers, data-source list and notes paragraph ---
$notes = $wb.Worksheets.Item("notes")

$notes.Range("A1").Value = "Notes"
$notes.Range("A5").Value = "Data Sources"
$notes.Range("A6").Value = "1. Netflix - searchlogistics.com; https://www.searchlogistics.com/learn/statistics/netflix-statistics/"
$notes.Range("A7").Value = "2. Prime - businessofapps.com; https://www.businessofapps.com/data/amazon-prime-video-statistics/"
$notes.Range("A8").Value = "3. Disney Plus - businessofapps.com; https://www.businessofapps.com/data/disney-plus-statistics/"
$notes.Range("A9").Value = "4. Hulu - businessofapps.com; https://www.businessofapps.com/data/hulu-statistics/"
$notes.Range("A10").Value = "5. Max - https://www.businessofapps.com/data/hbo-max-statistics/"
$notes.Range("A11").Value = "6. Paramount Plus - Wikipedia and associated links"
$notes.Range("A2").Value = "1. All figures are rounded to 1 decimal place.`n2. Users are used rather than subscriptions for Prime as Prime Video is a subsidiary of Amazon Prime.`n3. Paramount numbers are estimates configured from reports.`n4. For Hulu and Max, subscription figures are taken from Q4 of the year."

# widen column A and wrap / top-align the notes paragraph, matching the authored layout
$notes.Columns.Item(1).ColumnWidth = 77 - 5/7
$notes.Range("A2").WrapText = $true
$notes.Range("A2").VerticalAlignment = -4160
$notes.Rows.Item(2).RowHeight = 75

# --- restore the on-screen selections exactly as left by the author ---
$notes.Range("C9").Select() | Out-Null
$subs.Range("E8").Select() | Out-Null
